# Update "想去人数" (number of people wanting to go) values on the
# "展览" and "全部类型" sheets, mirroring the source data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 9
$ws1.Range("F5").Value = 899

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 9
$ws4.Range("F5").Value = 899
